$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 81 (shifts rows 81:199 down to 82:200)
$ws.Rows("81:81").Insert()

# Populate the new row 81 with the new weekly price record
$ws.Range("A81").Value = 11
$ws.Range("B81").Value = "Vega Monumental Concepción"
$ws.Range("C81").Value = "Bíobío"
$ws.Range("D81").Value = 44638
$ws.Range("E81").Value = 8
$ws.Range("F81").Value = 100114013
$ws.Range("G81").Value = "Zanahoria"
$ws.Range("H81").Value = "Sin especificar"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 250
$ws.Range("K81").Value = 7000
$ws.Range("L81").Value = 7500
$ws.Range("M81").Value = 7300
$ws.Range("N81").Value = "$/saco 20 kilos"
$ws.Range("O81").Value = "Chillán"
$ws.Range("P81").Value = 365
$ws.Range("Q81").Value = 20
$ws.Range("R81").Value = "Hortaliza"
